# Grupo4_LicoesAprendidas.pptx edit — "Exportando o slide em pdf para entrega"
#
# 1) Notes Master "datetimeFigureOut" field cached text: 2/5/2025 -> 2/9/2025
# 2) Slide 4, shape "Retângulo 1": widen the box (cx 1115512 -> 1292732 EMU)
# 3) Slide 4, shape "Retângulo 1": fix the wording
#    "Ouvir e ter a mente aberta" -> "Ouvir e manter a mente aberta"

$p = $ppt.ActivePresentation

# --- 1) Update the cached date/time field shown on the Notes pages ---
$nm = $p.NotesMaster
$dateField = $nm.HeadersFooters.DateAndTime
$dateField.Text = "2/9/2025"

# --- 2) & 3) Fix up the yellow/pink "Ouvir e ter a mente aberta" rectangle ---
$slide = $p.Slides.Item(4)
$shp = $slide.Shapes.Item("Retângulo 1")

$shp.Width = 101.7899212598425
$shp.TextFrame.TextRange.Text = "Ouvir e manter a mente aberta"
